$wb = $excel.ActiveWorkbook

# --- Sheet "packages": bump version/date in description ---
$pkg = $wb.Worksheets.Item("packages")
$pkg.Range("C2").Value = "Staging tables for raw data exports (v1.3.0, 2022-02-17)"

# --- Sheet "entities": add new entity row for cineasmappings ---
$ent = $wb.Worksheets.Item("entities")
$ent.Cells.Item(11, 1).Value = "cosasportal"
$ent.Cells.Item(11, 2).Value = "cineasmappings"
$ent.Cells.Item(11, 3).Value = "Cineas to HPO mappings"

# --- Sheet "attributes": add new attribute rows for cosasportal_cineasmappings ---
$attr = $wb.Worksheets.Item("attributes")

# row 111: value (id attribute)
$attr.Cells.Item(111, 1).Value = "cosasportal_cineasmappings"
$attr.Cells.Item(111, 2).Value = "value"
$attr.Cells.Item(111, 4).Value = $true
$attr.Cells.Item(111, 5).Value = $false
$attr.Cells.Item(111, 6).Value = $false
$attr.Cells.Item(111, 7).Value = $false
$attr.Cells.Item(111, 8).Value = "string"

# row 112: description
$attr.Cells.Item(112, 1).Value = "cosasportal_cineasmappings"
$attr.Cells.Item(112, 2).Value = "description"
$attr.Cells.Item(112, 4).Value = $false
$attr.Cells.Item(112, 5).Value = $false
$attr.Cells.Item(112, 6).Value = $true
$attr.Cells.Item(112, 7).Value = $false
$attr.Cells.Item(112, 8).Value = "string"

# row 113: codesystem
$attr.Cells.Item(113, 1).Value = "cosasportal_cineasmappings"
$attr.Cells.Item(113, 2).Value = "codesystem"
$attr.Cells.Item(113, 4).Value = $false
$attr.Cells.Item(113, 5).Value = $false
$attr.Cells.Item(113, 6).Value = $true
$attr.Cells.Item(113, 7).Value = $false
$attr.Cells.Item(113, 8).Value = "string"

# row 114: code
$attr.Cells.Item(114, 1).Value = "cosasportal_cineasmappings"
$attr.Cells.Item(114, 2).Value = "code"
$attr.Cells.Item(114, 4).Value = $false
$attr.Cells.Item(114, 5).Value = $false
$attr.Cells.Item(114, 6).Value = $true
$attr.Cells.Item(114, 7).Value = $false
$attr.Cells.Item(114, 8).Value = "string"

# row 115: hpo
$attr.Cells.Item(115, 1).Value = "cosasportal_cineasmappings"
$attr.Cells.Item(115, 2).Value = "hpo"
$attr.Cells.Item(115, 4).Value = $false
$attr.Cells.Item(115, 5).Value = $false
$attr.Cells.Item(115, 6).Value = $true
$attr.Cells.Item(115, 7).Value = $false
$attr.Cells.Item(115, 8).Value = "string"
